# Remove the depth, altitude, and elevation columns from the MIMS sheet
# since they are provided by the environmental package when relevant.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MIMS")

# Make MIMS the active sheet/tab (mirrors the workbook's activeTab + tabSelected move)
$ws.Activate() | Out-Null

# Select the depth/altitude/elevation columns (D:F) before removing them,
# matching the selection left behind after the deletion.
$ws.Range("D1:F1").EntireColumn.Select() | Out-Null

# Delete the depth, altitude and elevation columns entirely.
$ws.Range("D1:F1").EntireColumn.Delete() | Out-Null

Write-Host "Removed depth/altitude/elevation columns from MIMS sheet"
